$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceOne = 1 per call signature used below.
# Signature: Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#                     MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#                     Format, ReplaceWith, Replace)

# --- Paragraph 1 ("To be very honest, ...I reflect on all that I have learnt") ---
# Insert a comma after "and as such"
$d.Content.Find.Execute("and as such I had", $true, $false, $false, $false, $false, `
  $true, 1, $false, "and as such, I had", 2) | Out-Null

# Insert a comma after "However"
$d.Content.Find.Execute("knowledge to C++. However due to this project", $true, $false, $false, $false, $false, `
  $true, 1, $false, "knowledge to C++. However, due to this project", 2) | Out-Null

# Insert a comma after "due to this project"
$d.Content.Find.Execute("due to this project I was able", $true, $false, $false, $false, $false, `
  $true, 1, $false, "due to this project, I was able", 2) | Out-Null

# Insert a comma after "In this paper"
$d.Content.Find.Execute("CPP. In this paper I reflect", $true, $false, $false, $false, $false, `
  $true, 1, $false, "CPP. In this paper, I reflect", 2) | Out-Null

# --- Paragraph 2 ("First of all, ...without the use of hash maps") ---
# Fix "calculated" -> "calculate"
$d.Content.Find.Execute("calculated the distance", $true, $false, $false, $false, $false, `
  $true, 1, $false, "calculate the distance", 2) | Out-Null

# Insert a comma after "However, in this project"
$d.Content.Find.Execute("However, in this project I did not", $true, $false, $false, $false, $false, `
  $true, 1, $false, "However, in this project, I did not", 2) | Out-Null

# Insert a comma after "As such"
$d.Content.Find.Execute("Java project. As such I had to learn", $true, $false, $false, $false, $false, `
  $true, 1, $false, "Java project. As such, I had to learn", 2) | Out-Null

# --- Add two new reflective paragraphs after the second body paragraph ---
$secondPara = $d.Paragraphs(30)
$endOfSecond = $secondPara.Range
$endOfSecond.Collapse(0)
$endOfSecond.InsertParagraphAfter()

$thirdPara = $d.Paragraphs(31)
$thirdPara.Range.InsertAfter("Another challenge that I had was the difference in syntax between Java and C++(no matter how slight). There were times when it was an annoyance but with time I was able to overcome all these challenges.")

$endOfThird = $thirdPara.Range
$endOfThird.Collapse(0)
$endOfThird.InsertParagraphAfter()

$fourthPara = $d.Paragraphs(32)
$fourthPara.Range.InsertAfter("At the end of it all, I overcame all the challenges. I also found out in my research all that c++ is capable of(particularly game development) and I genuinely see a future for myself in that field with C++")

Write-Output "Paragraph count after edits: $($d.Paragraphs.Count)"
